# "fix: Commit correct results"
# The two benchmark blocks in the "results" sheet (B2:B20 and B23:B41) had
# been committed with the wrong measured numbers. Replace them with the
# correct values; the two line charts on the sheet read their series
# straight from these ranges, so their cached points are refreshed for free
# once the cells recalc.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

$block1 = @{
    2  = 51442187
    3  = 47049680
    4  = 43121466
    5  = 42999716
    6  = 43314530
    7  = 64611869
    8  = 68278414
    9  = 68768311
    10 = 80729080
    11 = 103031369
    12 = 110832062
    13 = 114464165
    14 = 116422450
    15 = 124861673
    16 = 249150229
    17 = 377067302
    18 = 353532543
    19 = 360349328
    20 = 373602684
}

$block2 = @{
    23 = 93036641
    24 = 94252564
    25 = 94446894
    26 = 94646841
    27 = 94847369
    28 = 99318001
    29 = 101212711
    30 = 102346371
    31 = 108746940
    32 = 120074375
    33 = 121873852
    34 = 122391716
    35 = 123048072
    36 = 128920942
    37 = 192796323
    38 = 265876588
    39 = 294861979
    40 = 304664490
    41 = 311390564
}

foreach ($r in $block1.Keys) {
    $ws.Cells.Item($r, 2).Value = $block1[$r]
}

foreach ($r in $block2.Keys) {
    $ws.Cells.Item($r, 2).Value = $block2[$r]
}

$wb.RefreshAll()

# Match the author's final selection state: the second block (B23:B41) is
# left selected on the data sheet.
$ws.Activate()
$ws.Range("B23:B41").Select()

$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Activate()
$ws2.Range("B23:B41").Select()

$ws.Activate()
